# Scheduled runner update: refresh market-board price columns (H-N) per leve row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(121, 8).Value = 497.6111
$ws.Cells.Item(121, 10).Value = 497.6111
$ws.Cells.Item(121, 12).Value = 1492.8333
$ws.Cells.Item(121, 14).Value = -4986.8333
$ws.Cells.Item(129, 8).Value = 1214.4894
$ws.Cells.Item(129, 10).Value = 1295.1395
$ws.Cells.Item(129, 12).Value = 3885.4185
$ws.Cells.Item(129, 14).Value = -13885.4185
$ws.Cells.Item(135, 8).Value = 765.8461
$ws.Cells.Item(135, 9).Value = 656.875
$ws.Cells.Item(135, 11).Value = 5911.875
$ws.Cells.Item(135, 13).Value = -3376.875
$ws.Cells.Item(137, 8).Value = 569489.1
$ws.Cells.Item(137, 9).Value = 1363246.2
$ws.Cells.Item(137, 10).Value = 2519.7551
$ws.Cells.Item(137, 11).Value = 4089738.6
$ws.Cells.Item(137, 12).Value = 7559.265299999999
$ws.Cells.Item(137, 13).Value = -4087188.6
$ws.Cells.Item(137, 14).Value = -12659.2653
$ws.Cells.Item(138, 8).Value = 2589.5217
$ws.Cells.Item(138, 9).Value = 1803.9333
$ws.Cells.Item(138, 10).Value = 4062.5
$ws.Cells.Item(138, 11).Value = 5411.7999
$ws.Cells.Item(138, 12).Value = 12187.5
$ws.Cells.Item(138, 13).Value = -271.7999
$ws.Cells.Item(138, 14).Value = -22467.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1961.6818
$ws.Cells.Item(61, 9).Value = 2026.8125
$ws.Cells.Item(61, 10).Value = 1788
$ws.Cells.Item(61, 11).Value = 2026.8125
$ws.Cells.Item(61, 12).Value = 1788
$ws.Cells.Item(61, 13).Value = -1814.8125
$ws.Cells.Item(61, 14).Value = -2212
$ws.Cells.Item(132, 8).Value = 2730.561
$ws.Cells.Item(132, 9).Value = 2060.9
$ws.Cells.Item(132, 10).Value = 4556.909
$ws.Cells.Item(132, 11).Value = 6182.700000000001
$ws.Cells.Item(132, 12).Value = 13670.727
$ws.Cells.Item(132, 13).Value = -3652.700000000001
$ws.Cells.Item(132, 14).Value = -18730.727
$ws.Cells.Item(136, 8).Value = 1961.6818
$ws.Cells.Item(136, 9).Value = 2026.8125
$ws.Cells.Item(136, 10).Value = 1788
$ws.Cells.Item(136, 11).Value = 6080.4375
$ws.Cells.Item(136, 12).Value = 5364
$ws.Cells.Item(136, 13).Value = -3530.4375
$ws.Cells.Item(136, 14).Value = -10464

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(75, 8).Value = 6876.909
$ws.Cells.Item(75, 9).Value = 4650.1665
$ws.Cells.Item(75, 10).Value = 9549
$ws.Cells.Item(75, 11).Value = 4650.1665
$ws.Cells.Item(75, 12).Value = 9549
$ws.Cells.Item(75, 13).Value = -3714.1665
$ws.Cells.Item(75, 14).Value = -11421
$ws.Cells.Item(78, 8).Value = 6876.909
$ws.Cells.Item(78, 9).Value = 4650.1665
$ws.Cells.Item(78, 10).Value = 9549
$ws.Cells.Item(78, 11).Value = 13950.4995
$ws.Cells.Item(78, 12).Value = 28647
$ws.Cells.Item(78, 13).Value = -9270.499500000002
$ws.Cells.Item(78, 14).Value = -38007
$ws.Cells.Item(99, 8).Value = 2135.5264
$ws.Cells.Item(99, 9).Value = 1128.2142
$ws.Cells.Item(99, 11).Value = 1128.2142
$ws.Cells.Item(99, 13).Value = 369.7858000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 5292865.5
$ws.Cells.Item(16, 9).Value = 10102447
$ws.Cells.Item(16, 10).Value = 2326.2
$ws.Cells.Item(16, 11).Value = 10102447
$ws.Cells.Item(16, 12).Value = 2326.2
$ws.Cells.Item(16, 13).Value = -10102160
$ws.Cells.Item(16, 14).Value = -2900.2
$ws.Cells.Item(31, 8).Value = 195168.78
$ws.Cells.Item(31, 9).Value = 483488
$ws.Cells.Item(31, 10).Value = 2955.976
$ws.Cells.Item(31, 11).Value = 483488
$ws.Cells.Item(31, 12).Value = 2955.976
$ws.Cells.Item(31, 13).Value = -483193
$ws.Cells.Item(31, 14).Value = -3545.976
$ws.Cells.Item(34, 8).Value = 195168.78
$ws.Cells.Item(34, 9).Value = 483488
$ws.Cells.Item(34, 10).Value = 2955.976
$ws.Cells.Item(34, 11).Value = 483488
$ws.Cells.Item(34, 12).Value = 2955.976
$ws.Cells.Item(34, 13).Value = -483286
$ws.Cells.Item(34, 14).Value = -3359.976
$ws.Cells.Item(41, 8).Value = 50000
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 10).Value = 50000
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 13).Value = 50000
$ws.Cells.Item(41, 14).Value = -50856
$ws.Cells.Item(68, 8).Value = 74999.5
$ws.Cells.Item(68, 10).Value = 74999.5
$ws.Cells.Item(68, 12).Value = 74999.5
$ws.Cells.Item(68, 14).Value = -76497.5
$ws.Cells.Item(71, 8).Value = 74999.5
$ws.Cells.Item(71, 10).Value = 74999.5
$ws.Cells.Item(71, 12).Value = 224998.5
$ws.Cells.Item(71, 14).Value = -232486.5
$ws.Cells.Item(87, 8).Value = 22565
$ws.Cells.Item(87, 10).Value = 22565
$ws.Cells.Item(87, 12).Value = 22565
$ws.Cells.Item(87, 14).Value = -24937
$ws.Cells.Item(90, 8).Value = 22565
$ws.Cells.Item(90, 10).Value = 22565
$ws.Cells.Item(90, 12).Value = 67695
$ws.Cells.Item(90, 14).Value = -79551
$ws.Cells.Item(113, 8).Value = 5292865.5
$ws.Cells.Item(113, 9).Value = 10102447
$ws.Cells.Item(113, 10).Value = 2326.2
$ws.Cells.Item(113, 11).Value = 10102447
$ws.Cells.Item(113, 12).Value = 2326.2
$ws.Cells.Item(113, 13).Value = -10100277
$ws.Cells.Item(113, 14).Value = -6666.2
$ws.Cells.Item(141, 8).Value = 27783.334
$ws.Cells.Item(141, 10).Value = 27783.334
$ws.Cells.Item(141, 12).Value = 27783.334
$ws.Cells.Item(141, 14).Value = -38143.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 2100909
$ws.Cells.Item(2, 10).Value = 2976266.5
$ws.Cells.Item(2, 12).Value = 17857599
$ws.Cells.Item(2, 14).Value = -17857825
$ws.Cells.Item(5, 8).Value = 557726.75
$ws.Cells.Item(5, 9).Value = 1668
$ws.Cells.Item(5, 10).Value = 704058
$ws.Cells.Item(5, 11).Value = 5004
$ws.Cells.Item(5, 12).Value = 2112174
$ws.Cells.Item(5, 13).Value = -4892
$ws.Cells.Item(5, 14).Value = -2112398
$ws.Cells.Item(56, 8).Value = 4854.2856
$ws.Cells.Item(56, 9).Value = 4854.2856
$ws.Cells.Item(56, 11).Value = 4854.2856
$ws.Cells.Item(56, 13).Value = -4324.2856
$ws.Cells.Item(59, 8).Value = 3984.1667
$ws.Cells.Item(59, 9).Value = 2968.3333
$ws.Cells.Item(59, 11).Value = 8904.999899999999
$ws.Cells.Item(59, 13).Value = -8364.999899999999
$ws.Cells.Item(68, 8).Value = 3661.3408
$ws.Cells.Item(68, 9).Value = 1003.0303
$ws.Cells.Item(68, 10).Value = 11636.272
$ws.Cells.Item(68, 11).Value = 3009.0909
$ws.Cells.Item(68, 12).Value = 34908.81600000001
$ws.Cells.Item(68, 13).Value = -2198.0909
$ws.Cells.Item(68, 14).Value = -36530.81600000001
$ws.Cells.Item(71, 8).Value = 3661.3408
$ws.Cells.Item(71, 9).Value = 1003.0303
$ws.Cells.Item(71, 10).Value = 11636.272
$ws.Cells.Item(71, 11).Value = 9027.2727
$ws.Cells.Item(71, 12).Value = 104726.448
$ws.Cells.Item(71, 13).Value = -4971.2727
$ws.Cells.Item(71, 14).Value = -112838.448
$ws.Cells.Item(126, 8).Value = 4016
$ws.Cells.Item(126, 9).Value = 2520
$ws.Cells.Item(126, 10).Value = 10000
$ws.Cells.Item(126, 11).Value = 7560
$ws.Cells.Item(126, 12).Value = 30000
$ws.Cells.Item(126, 13).Value = -2620
$ws.Cells.Item(126, 14).Value = -39880
$ws.Cells.Item(135, 8).Value = 557726.75
$ws.Cells.Item(135, 9).Value = 1668
$ws.Cells.Item(135, 10).Value = 704058
$ws.Cells.Item(135, 11).Value = 15012
$ws.Cells.Item(135, 12).Value = 6336522
$ws.Cells.Item(135, 13).Value = -12477
$ws.Cells.Item(135, 14).Value = -6341592

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 39666.668
$ws.Cells.Item(15, 10).Value = 39666.668
$ws.Cells.Item(15, 12).Value = 39666.668
$ws.Cells.Item(15, 14).Value = -40242.668
$ws.Cells.Item(81, 8).Value = 39666.668
$ws.Cells.Item(81, 10).Value = 39666.668
$ws.Cells.Item(81, 12).Value = 39666.668
$ws.Cells.Item(81, 14).Value = -41662.668
$ws.Cells.Item(84, 8).Value = 39666.668
$ws.Cells.Item(84, 10).Value = 39666.668
$ws.Cells.Item(84, 12).Value = 119000.004
$ws.Cells.Item(84, 14).Value = -128984.004
$ws.Cells.Item(135, 8).Value = 60000
$ws.Cells.Item(135, 10).Value = 60000
$ws.Cells.Item(135, 12).Value = 60000
$ws.Cells.Item(135, 14).Value = -70140

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 5293269
$ws.Cells.Item(93, 9).Value = 7409423
$ws.Cells.Item(93, 11).Value = 7409423
$ws.Cells.Item(93, 13).Value = -7408175
$ws.Cells.Item(129, 8).Value = 41686.668
$ws.Cells.Item(129, 10).Value = 41686.668
$ws.Cells.Item(129, 12).Value = 41686.668
$ws.Cells.Item(129, 14).Value = -51686.668
$ws.Cells.Item(141, 8).Value = 32275
$ws.Cells.Item(141, 10).Value = 32275
$ws.Cells.Item(141, 12).Value = 32275
$ws.Cells.Item(141, 14).Value = -42635

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1333209.9
$ws.Cells.Item(126, 9).Value = 2076.4
$ws.Cells.Item(126, 10).Value = 3551765.8
$ws.Cells.Item(126, 11).Value = 6229.200000000001
$ws.Cells.Item(126, 12).Value = 10655297.4
$ws.Cells.Item(126, 13).Value = -3759.200000000001
$ws.Cells.Item(126, 14).Value = -10660237.4
$ws.Cells.Item(136, 8).Value = 2301.8
$ws.Cells.Item(136, 9).Value = 928.9032
$ws.Cells.Item(136, 10).Value = 4541.7896
$ws.Cells.Item(136, 11).Value = 2786.7096
$ws.Cells.Item(136, 12).Value = 13625.3688
$ws.Cells.Item(136, 13).Value = -236.7096000000001
$ws.Cells.Item(136, 14).Value = -18725.3688
